{"js": "const body = context.document.body;\n\nconst replacements = [\n  [\"2025-05-15 Thursday\", \"2025-05-16 Friday\"],\n  [\"93\u00f78=\", \"61\u00f74=\"],\n  [\"71\u00f74=\", \"58\u00f77=\"],\n  [\"84\u00f73=\", \"49\u00f75=\"],\n  [\"20\u00f72=\", \"35\u00f75=\"],\n  [\"26\u00f76=\", \"51\u00f77=\"],\n  [\"18\u00f73=\", \"22\u00f79=\"],\n  [\"51\u00f73=\", \"69\u00f77=\"],\n  [\"12\u00f78=\", \"18\u00f77=\"],\n  [\"79\u00f79=\", \"99\u00f78=\"],\n  [\"54\u00f73=\", \"31\u00f73=\"],\n  [\"29\u00f76=\", \"15\u00f73=\"],\n  [\"42\u00f73=\", \"76\u00f78=\"],\n  [\"53\u00f74=\", \"84\u00f72=\"],\n  [\"47\u00f74=\", \"87\u00f77=\"],\n  [\"88\u00f78=\", \"23\u00f77=\"],\n  [\"62\u00f78=\", \"31\u00f76=\"],\n  [\"55\u00f75=\", \"37\u00f75=\"],\n  [\"93\u00f77=\", \"59\u00f76=\"],\n  [\"55\u00f73=\", \"61\u00f76=\"],\n  [\"96\u00f72=\", \"90\u00f76=\"],\n  [\"37\u00f78=\", \"98\u00f78=\"],\n  [\"98\u00f77=\", \"65\u00f77=\"],\n  [\"14\u00f73=\", \"37\u00f75=\"],\n  [\"42\u00f76=\", \"63\u00f72=\"],\n  [\"87\u00f76=\", \"88\u00f73=\"],\n];\n\nfor (const [findText, replaceText] of replacements) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    continue;\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replaceText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-AllText($doc, $findText, $replaceText) {\n    $find = $doc.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Execute([ref]$findText, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$replaceText, 2) | Out-Null\n}\n\nReplace-AllText $d '2025-05-15 Thursday' '2025-05-16 Friday'\nReplace-AllText $d '93\u00f78=' '61\u00f74='\nReplace-AllText $d '71\u00f74=' '58\u00f77='\nReplace-AllText $d '84\u00f73=' '49\u00f75='\nReplace-AllText $d '20\u00f72=' '35\u00f75='\nReplace-AllText $d '26\u00f76=' '51\u00f77='\nReplace-AllText $d '18\u00f73=' '22\u00f79='\nReplace-AllText $d '51\u00f73=' '69\u00f77='\nReplace-AllText $d '12\u00f78=' '18\u00f77='\nReplace-AllText $d '79\u00f79=' '99\u00f78='\nReplace-AllText $d '54\u00f73=' '31\u00f73='\nReplace-AllText $d '29\u00f76=' '15\u00f73='\nReplace-AllText $d '42\u00f73=' '76\u00f78='\nReplace-AllText $d '53\u00f74=' '84\u00f72='\nReplace-AllText $d '47\u00f74=' '87\u00f77='\nReplace-AllText $d '88\u00f78=' '23\u00f77='\nReplace-AllText $d '62\u00f78=' '31\u00f76='\nReplace-AllText $d '55\u00f75=' '37\u00f75='\nReplace-AllText $d '93\u00f77=' '59\u00f76='\nReplace-AllText $d '55\u00f73=' '61\u00f76='\nReplace-AllText $d '96\u00f72=' '90\u00f76='\nReplace-AllText $d '37\u00f78=' '98\u00f78='\nReplace-AllText $d '98\u00f77=' '65\u00f77='\nReplace-AllText $d '14\u00f73=' '37\u00f75='\nReplace-AllText $d '42\u00f76=' '63\u00f72='\nReplace-AllText $d '87\u00f76=' '88\u00f73='\n"}
